$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1659.8
$ws.Range("J32").Value = 1600
$ws.Range("L32").Value = 1600
$ws.Range("N32").Value = -2252

$ws.Range("H38").Value = 3307.7778
$ws.Range("I38").Value = 3223
$ws.Range("K38").Value = 9669
$ws.Range("M38").Value = -9297

$ws.Range("H100").Value = 9093.789000000001
$ws.Range("I100").Value = 1713.4286
$ws.Range("K100").Value = 1713.4286
$ws.Range("M100").Value = -1172.4286

$ws.Range("H106").Value = 4950
$ws.Range("I106").Value = 4950
$ws.Range("K106").Value = 4950
$ws.Range("M106").Value = -4319

$ws.Range("H132").Value = 2065.9688
$ws.Range("I132").Value = 2263.2222
$ws.Range("K132").Value = 6789.6666
$ws.Range("M132").Value = -4259.6666

$ws.Range("H137").Value = 4461.391
$ws.Range("I137").Value = 4362.476
$ws.Range("K137").Value = 13087.428
$ws.Range("M137").Value = -10537.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2846.3635
$ws.Range("I61").Value = 2431
$ws.Range("K61").Value = 2431
$ws.Range("M61").Value = -2219

$ws.Range("H132").Value = 2550.5557
$ws.Range("I132").Value = 2550.5557
$ws.Range("K132").Value = 7651.6671
$ws.Range("M132").Value = -5121.6671

$ws.Range("H136").Value = 2846.3635
$ws.Range("I136").Value = 2431
$ws.Range("K136").Value = 7293
$ws.Range("M136").Value = -4743

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2695
$ws.Range("I99").Value = 3632
$ws.Range("J99").Value = 2070.3333
$ws.Range("K99").Value = 3632
$ws.Range("L99").Value = 2070.3333
$ws.Range("M99").Value = -2134
$ws.Range("N99").Value = -5066.3333

$ws.Range("H134").Value = 44743.957
$ws.Range("I134").Value = 3211.087
$ws.Range("K134").Value = 9633.261
$ws.Range("M134").Value = -7098.261

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3559.4
$ws.Range("I16").Value = 2865.6667
$ws.Range("J16").Value = 4600
$ws.Range("K16").Value = 2865.6667
$ws.Range("L16").Value = 4600
$ws.Range("M16").Value = -2578.6667
$ws.Range("N16").Value = -5174

$ws.Range("H58").Value = 6476.525
$ws.Range("I58").Value = 6316.7407
$ws.Range("K58").Value = 6316.7407
$ws.Range("M58").Value = -6113.7407

$ws.Range("H113").Value = 3559.4
$ws.Range("I113").Value = 2865.6667
$ws.Range("J113").Value = 4600
$ws.Range("K113").Value = 2865.6667
$ws.Range("L113").Value = 4600
$ws.Range("M113").Value = -695.6667000000002
$ws.Range("N113").Value = -8940

$ws.Range("H116").Value = 64221
$ws.Range("J116").Value = 64221
$ws.Range("L116").Value = 64221
$ws.Range("N116").Value = -73399

$ws.Range("H122").Value = 4011.7334
$ws.Range("I122").Value = 4298.4
$ws.Range("J122").Value = 3868.4
$ws.Range("K122").Value = 12895.2
$ws.Range("L122").Value = 11605.2
$ws.Range("M122").Value = -10445.2
$ws.Range("N122").Value = -16505.2

$ws.Range("H132").Value = 1598.04
$ws.Range("I132").Value = 1620.3182
$ws.Range("J132").Value = 1434.6666
$ws.Range("K132").Value = 4860.9546
$ws.Range("L132").Value = 4303.9998
$ws.Range("M132").Value = -2330.9546
$ws.Range("N132").Value = -9363.9998

$ws.Range("H134").Value = 403379.8
$ws.Range("I134").Value = 3520.625
$ws.Range("K134").Value = 10561.875
$ws.Range("M134").Value = -8026.875

$ws.Range("H136").Value = 6476.525
$ws.Range("I136").Value = 6316.7407
$ws.Range("K136").Value = 18950.2221
$ws.Range("M136").Value = -16400.2221

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 201.46666
$ws.Range("I12").Value = 200
$ws.Range("J12").Value = 201.57143
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 604.71429
$ws.Range("M12").Value = -427
$ws.Range("N12").Value = -950.71429

$ws.Range("H38").Value = 31
$ws.Range("J38").Value = 14
$ws.Range("L38").Value = 42
$ws.Range("N38").Value = -736

$ws.Range("H104").Value = 1869
$ws.Range("I104").Value = 1813
$ws.Range("J104").Value = 1925
$ws.Range("K104").Value = 5439
$ws.Range("L104").Value = 5775
$ws.Range("M104").Value = -2818
$ws.Range("N104").Value = -11017

$ws.Range("H114").Value = 286.55554
$ws.Range("I114").Value = 224.5
$ws.Range("J114").Value = 336.2
$ws.Range("K114").Value = 673.5
$ws.Range("L114").Value = 1008.6
$ws.Range("M114").Value = 2580.5
$ws.Range("N114").Value = -7516.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2275.4375
$ws.Range("I102").Value = 1113.25
$ws.Range("J102").Value = 3437.625
$ws.Range("K102").Value = 1113.25
$ws.Range("L102").Value = 3437.625
$ws.Range("M102").Value = 508.75
$ws.Range("N102").Value = -6681.625

$ws.Range("H126").Value = 38464890
$ws.Range("I126").Value = 66669540
$ws.Range("K126").Value = 200008620
$ws.Range("M126").Value = -200006150

$ws.Range("H129").Value = 76000
$ws.Range("J129").Value = 76000
$ws.Range("L129").Value = 76000
$ws.Range("N129").Value = -86000

$ws.Range("H132").Value = 46804.207
$ws.Range("I132").Value = 4840.706
$ws.Range("J132").Value = 148715.58
$ws.Range("K132").Value = 14522.118
$ws.Range("L132").Value = 446146.74
$ws.Range("M132").Value = -11992.118
$ws.Range("N132").Value = -451206.74

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 75000
$ws.Range("J69").Value = 75000
$ws.Range("L69").Value = 75000
$ws.Range("N69").Value = -76622

$ws.Range("H72").Value = 75000
$ws.Range("J72").Value = 75000
$ws.Range("L72").Value = 225000
$ws.Range("N72").Value = -233112

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1711.6471
$ws.Range("I81").Value = 1647.6154
$ws.Range("K81").Value = 3295.2308
$ws.Range("M81").Value = -2234.2308

$ws.Range("H84").Value = 1711.6471
$ws.Range("I84").Value = 1647.6154
$ws.Range("K84").Value = 16476.154
$ws.Range("M84").Value = -11172.154

$ws.Range("H105").Value = 85615
$ws.Range("J105").Value = 85615
$ws.Range("L105").Value = 85615
$ws.Range("N105").Value = -92603

$ws.Range("H107").Value = 598.94116
$ws.Range("I107").Value = 636.0741
$ws.Range("J107").Value = 455.7143
$ws.Range("K107").Value = 1908.2223
$ws.Range("L107").Value = 1367.1429
$ws.Range("M107").Value = 11.77769999999987
$ws.Range("N107").Value = -5207.1429

$ws.Range("H132").Value = 27687.95
$ws.Range("I132").Value = 2447.25
$ws.Range("J132").Value = 86582.914
$ws.Range("K132").Value = 7341.75
$ws.Range("L132").Value = 259748.742
$ws.Range("M132").Value = -4811.75
$ws.Range("N132").Value = -264808.742

$ws.Range("H136").Value = 9346045
$ws.Range("I136").Value = 11459742
$ws.Range("J136").Value = 287344
$ws.Range("K136").Value = 34379226
$ws.Range("L136").Value = 862032
$ws.Range("M136").Value = -34376676
$ws.Range("N136").Value = -867132
